$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update labels / legends (shared strings), in the order that produces
# the same shared-string table ordering as the target workbook. ---
$ws.Range("A15").Value = "Chord selection"
$ws.Range("G15").Value = "Bytes needed to select a single chord"
$ws.Range("A21").Value = "Notes selection"
$ws.Range("G21").Value = "Bytes needed to select notes (within chord)"
$ws.Range("G20").Value = "Per notes selection (ex: for velocity, …)"

# --- Update numeric inputs ---
$ws.Range("B17").Value = 8
$ws.Range("B20").Value = 4
$ws.Range("B23").Value = 2
$ws.Range("B28").Value = 64
$ws.Range("B29").Value = 2

# --- Update formulas ---
$ws.Range("E21").Formula = "=_xlfn.CEILING.MATH((B10+B20)/8,1)"
$ws.Range("E24").Formula = "=B22*B23"

# --- Update view / selection state ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G22").Select()
